$d = $word.ActiveDocument

# 1. Table column width: 2795 dxa -> 2869 dxa (1 pt = 20 dxa)
$t = $d.Tables.Item(1)
$col = $t.Columns.Item(3)
$col.Width = 2869 / 20

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 2 years block
Replace-Text "2.59%" "2.71%"
Replace-Text "2.6% (2.41% to 2.78%)" "2.73% (2.54% to 2.92%)"
Replace-Text "1 (0.93 to 1.08)" "1.01 (0.94 to 1.08)"
Replace-Text "0% (-0.19% to 0.19%)" "0.02% (-0.17% to 0.21%)"
Replace-Text "0 (-0.1 to 0.1)" "-0.01 (-0.11 to 0.1)"
Replace-Text "-0.14 (-0.21 to -0.07)" "-0.15 (-0.24 to -0.06)"
Replace-Text "0.91 (0.89 to 0.92)" "0.9 (0.89 to 0.92)"

# 5 years block
Replace-Text "4.39%" "4.72%"
Replace-Text "4.49% (4.24% to 4.74%)" "4.76% (4.51% to 5.02%)"
Replace-Text "1.02 (0.97 to 1.08)" "1.01 (0.96 to 1.07)"
Replace-Text "0.1% (-0.15% to 0.35%)" "0.05% (-0.21% to 0.3%)"
Replace-Text "-0.02 (-0.1 to 0.06)" "-0.02 (-0.09 to 0.06)"
Replace-Text "-0.16 (-0.22 to -0.1)" "-0.16 (-0.21 to -0.1)"
Replace-Text "0.89 (0.88 to 0.9)" "0.88 (0.87 to 0.89)"
Replace-Text "0.04 (0.03 to 0.04)" "0.04 (0.04 to 0.04)"

Write-Output "done"
